# Update "Latest HO Xliff Generate Date" / "Correspond Handoff/Handback Datetime"
# timestamps for the ed2351ae-a01f-4cb1-a7c6-0a3c024b7d62 row across the
# Overview, zh-cn, and de-de sheets, reflecting a newer handback report run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview!G4 - Latest HO Xliff Generate Date
$overview.Range("G4").Value = "2016-09-04 12:49:07"

# zh-cn!H4 - Correspond Handoff Datetime
$zhcn.Range("H4").Value = "2016-09-04 12:48:58"

# zh-cn!K4 - Correspond Handback DateTime
$zhcn.Range("K4").Value = "2016-09-04 12:49:28"

# de-de!H4 - Correspond Handoff Datetime (shared text with Overview!G4)
$dede.Range("H4").Value = "2016-09-04 12:49:07"

# de-de!K4 - Correspond Handback DateTime
$dede.Range("K4").Value = "2016-09-04 12:49:35"
